$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append: date serial, B, C, D
$data = @(
    @(44330, 0, 6, 38.75217981011431),
    @(44331, 4, 8, 51.66957308015243),
    @(44332, 0, 7, 45.21087644513337),
    @(44333, 0, 6, 38.75217981011431),
    @(44334, 0, 5, 32.29348317509527),
    @(44335, 0, 5, 32.29348317509527),
    @(44336, 3, 7, 45.21087644513337),
    @(44337, 1, 8, 51.66957308015243),
    @(44338, 4, 8, 51.66957308015243),
    @(44339, 0, 8, 51.66957308015243),
    @(44340, 2, 10, 64.58696635019054),
    @(44341, 0, 10, 64.58696635019054),
    @(44342, 0, 10, 64.58696635019054),
    @(44343, 0, 7, 45.21087644513337)
)

$startRow = 256
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $cellA = $ws.Cells.Item($row, 1)
    $ws.Cells.Item(255, 1).Copy($cellA)
    $cellA.Value = $vals[0]

    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}
